# Generate Report for Handoff
# Replace the GUID-based file identifier c76d78d2-8e68-4bb0-945d-74475564cb3b
# with the freshly generated 919668ce-0deb-44e5-82e5-c704e015e9a7, update the
# xliff content hashes, and bump the handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "c76d78d2-8e68-4bb0-945d-74475564cb3b"
$newGuid = "919668ce-0deb-44e5-82e5-c704e015e9a7"

$oldZhXlf = "$oldGuid.e4f0a6fd0184921f993ca6c79a2b980c74588050.zh-cn.xlf"
$newZhXlf = "$newGuid.3df374e71a09553796dbe5c4c0a9960c02d4e655.zh-cn.xlf"

$oldDeXlf = "$oldGuid.e4f0a6fd0184921f993ca6c79a2b980c74588050.de-de.xlf"
$newDeXlf = "$newGuid.3df374e71a09553796dbe5c4c0a9960c02d4e655.de-de.xlf"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-27 19:04:35"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-27 19:04:30"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-08-27 19:04:35"
